# "break out stock.yaml completed"
# Target sheet is the "week" tab (sheet2.xml), which currently spans A1:I82.
# 1) Rows 72-82 had their D column (bsecode) stored as text; convert to numeric.
# 2) Append a fresh block of 11 rows (83-93) that duplicates rows 72-82 (new
#    scrape of the same stocks), keeping D/H/I as text like the original
#    pre-conversion rows, but with updated volume (G) and timestamp (I).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week")

# Step 1: duplicate rows 72:82 into rows 83:93 *before* changing any types,
# so the copy inherits the original text-typed D/H/I cells verbatim.
$ws.Range("A72:I82").Copy()
$ws.Range("A83").PasteSpecial(-4104)
$excel.CutCopyMode = 0

# Step 2: convert D72:D82 (bsecode) from text to real numbers.
$ws.Range("D72").Value = 539523
$ws.Range("D73").Value = 500257
$ws.Range("D74").Value = 524715
$ws.Range("D75").Value = 532321
$ws.Range("D76").Value = 542830
$ws.Range("D77").Value = 500295
$ws.Range("D78").Value = 540222
$ws.Range("D79").Value = 513599
$ws.Range("D80").Value = 500103
$ws.Range("D81").Value = 526371
$ws.Range("D82").Value = 500113

# Step 3: fix up the volume (G) figures that changed between scrapes.
$ws.Range("G83").Value = 272426
$ws.Range("G90").Value = 6388899
$ws.Range("G93").Value = 72583869

# Step 4: stamp the new rows with the later scrape timestamp.
foreach ($r in 83..93) {
    $ws.Cells.Item($r, 9).Value = "30/06/2024 18:35:16"
}
